$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column holds text-formatted price strings (e.g. "1.00", "95.188.10").
# Excel auto-converts numeric-looking text typed into .Value, which would
# silently drop significant trailing zeros / reparse dotted price strings
# as numbers. Force text entry via a transient Text number format, then
# restore the cell to the workbooks default (unstyled) "Normal" style so
# we do not leave a stray number-format behind on the cell.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "95.188.10"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.97%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.608.09"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -2.73%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "2.28"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +19.64%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "225.99"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.23%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "634.46"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -3.44%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.411"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.27%  "

$ws.Range("E9").Value = "  +1.82%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "3.610.97"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.58%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "46.95"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +5.59%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.206"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.83%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000289"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.36%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.46"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.52%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "4.279.38"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.75%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "94.842.87"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.00%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "8.77"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.32%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.618.97"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.26%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "19.27"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.60%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.62"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.73%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.511"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "510.13"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.82%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.23"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -5.72%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.239"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +22.26%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "113.05"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +11.02%  "

$ws.Range("E27").Value = "  -5.64%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.71"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.69%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "12.58"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -6.54%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "12.61"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.68%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.17%  "

$ws.Range("E33").Value = "  +0.24%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.178"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -6.35%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.76"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -6.22%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "31.72"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.06%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.583"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.25%  "

$ws.Range("E38").Value = "  -0.03%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "593.14"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -8.87%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "8.28"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -6.63%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.80"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.480"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +7.46%  "

$ws.Range("E43").Value = "  -1.97%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "39.32"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.42%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0477"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.32%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.92"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -6.30%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.913"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -5.10%  "

$ws.Range("E48").Value = "  -0.86%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.53"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.19"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.51%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "53.81"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
